# "prepped for presentation and reorganized"
# - Row 18 (2/23 entry) gets a running total note appended to its date cell.
# - A new row 19 is added for 2/24 (45 mins), reorganizing the app and
#   taking pictures for the presentation.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the 2/23 date cell to note the running total.
$ws.Cells.Item(18, 1).Value = "2/23, 2 hours (4 total)"

# Bring row 19 into existence with the same look (borders/number format/
# wrap) as the row above it, then fill in the new log entry.
$ws.Range("A18:B18").Copy()
$ws.Range("A19:B19").PasteSpecial(-4122)

$ws.Cells.Item(19, 1).Value = "2/24, 45 mins"
$ws.Cells.Item(19, 2).Value = "Worked on reorganizing app, taking pics for presentation"

# Mirror the author's final selection after typing the new row.
$ws.Range("B20").Select()
